$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, derived from the target diff.
$changes = @{
    "D2" = "53.982.80"
    "E2" = "  -11.03%  "
    "D3" = "2.334.30"
    "E3" = "  -19.69%  "
    "E4" = "  +0.07%  "
    "D5" = "444.19"
    "E5" = "  -16.05%  "
    "D6" = "126.78"
    "E6" = "  -12.07%  "
    "E7" = "  -0.22%  "
    "D8" = "0.476"
    "E8" = "  -14.58%  "
    "D9" = "2.345.19"
    "E9" = "  -19.58%  "
    "D10" = "5.38"
    "E10" = "  -11.11%  "
    "D11" = "0.0918"
    "E11" = "  -15.49%  "
    "D12" = "0.309"
    "E12" = "  -14.83%  "
    "E13" = "  -3.26%  "
    "D14" = "2.741.92"
    "E14" = "  -19.75%  "
    "D15" = "53.997.30"
    "E15" = "  -10.96%  "
    "D16" = "18.85"
    "E16" = "  -17.51%  "
    "E17" = "  -14.48%  "
    "D18" = "2.359.27"
    "E18" = "  -19.02%  "
    "D19" = "3.95"
    "E19" = "  -21.74%  "
    "D20" = "298.05"
    "E20" = "  -18.00%  "
    "D21" = "9.16"
    "E21" = "  -21.95%  "
    "E22" = "  -0.20%  "
    "D23" = "5.58"
    "E23" = "  -1.72%  "
    "D24" = "5.39"
    "E24" = "  -18.87%  "
    "D25" = "55.68"
    "E25" = "  -14.26%  "
    "D26" = "0.999"
    "E26" = "  +0.11%  "
    "D27" = "0.152"
    "E27" = "  -16.20%  "
    "D28" = "0.367"
    "E28" = "  -19.52%  "
    "D29" = "6.95"
    "E29" = "  -11.80%  "
    "E30" = "  -0.29%  "
    "D31" = "0.0₃0702"
    "E31" = "  -18.69%  "
    "D32" = "146.40"
    "E32" = "  -3.93%  "
    "D33" = "17.26"
    "E33" = "  -12.79%  "
    "D34" = "1.35"
    "E34" = "  -19.64%  "
    "D35" = "4.66"
    "E35" = "  -16.75%  "
    "D36" = "3.55"
    "E36" = "  -19.17%  "
    "D37" = "0.832"
    "E37" = "  -17.73%  "
    "E38" = "  -17.03%  "
    "D39" = "33.34"
    "E39" = "  -11.64%  "
    "D40" = "0.996"
    "E40" = "  -0.11%  "
    "E41" = "  -0.60%  "
    "D42" = "3.13"
    "E42" = "  -16.11%  "
    "D43" = "1.938.19"
    "E43" = "  -15.76%  "
    "D44" = "1.21"
    "E44" = "  -18.71%  "
    "D45" = "0.0494"
    "E45" = "  -15.20%  "
    "D46" = "0.527"
    "E46" = "  -18.87%  "
    "D47" = "0.0210"
    "E47" = "  -11.74%  "
    "D48" = "0.0833"
    "E48" = "  -10.10%  "
    "B49" = "RenderToken"
    "C49" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D49" = "4.04"
    "E49" = "  -19.78%  "
    "B50" = "EnergySwap"
    "C50" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D50" = "15.83"
    "E50" = "  -23.06%  "
    "E51" = "  -3.62%  "
}

# The "Price" column (D) sometimes holds values that look like plain
# numbers (e.g. "0.476", "444.19"). Excel normally auto-converts such
# text into a numeric value when assigned through .Value. Since the
# source data must stay as literal text (matching the original sheet,
# which stores every cell as a string), temporarily mark the Price
# column as Text-formatted while writing, then restore the original
# "Normal" style so no visible formatting change remains.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}

$priceRange.Style = "Normal"

Write-Host "Applied $($changes.Count) cell updates"
